# Commit: "Added more changes and a git hook."
#
# Semantic edit performed:
#   1. The (first) paragraph reading "This is a sub-sub test" has its single
#      run split into four runs - "This is a sub-", "sub", " ", "change" -
#      i.e. the text becomes "This is a sub-sub change", authored run-by-run
#      the way Word leaves things behind when a user edits/accepts changes
#      incrementally.
#   2. The document's "_GoBack" bookmark (which marks the most recent edit
#      location) moves from the end of the "Added this new change" paragraph
#      to the end of the newly-edited "This is a sub-sub change" paragraph.

$d = $word.ActiveDocument

function Get-ParagraphPPr($para) {
    $xml = $para.Range.WordOpenXML
    if ($xml -match '(?s)<w:pPr>.*?</w:pPr>') {
        return $matches[0]
    }
    return ""
}

function Get-ParagraphText($para) {
    return $para.Range.Text.TrimEnd([char]13, [char]7)
}

function Get-XmlEscapedRun($text) {
    $escaped = $text.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    if ($text -match '^\s' -or $text -match '\s$' -or $text -eq "") {
        return "<w:r><w:t xml:space=`"preserve`">$escaped</w:t></w:r>"
    }
    return "<w:r><w:t>$escaped</w:t></w:r>"
}

function New-PackageXml($bodyXml) {
    return '<?xml version="1.0" standalone="yes"?>' + `
        '<?mso-application progid="Word.Document"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        "<w:body>$bodyXml</w:body>" + `
        '</w:document>' + `
        '</pkg:xmlData>' + `
        '</pkg:part>' + `
        '</pkg:package>'
}

# --- locate the paragraphs we need to touch -------------------------------
$subSubIdx = -1
$addedIdx = -1

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $t = Get-ParagraphText $p

    if ($subSubIdx -eq -1 -and $t -eq "This is a sub-sub test") {
        $subSubIdx = $i
    }
    if ($t -eq "Added this new change") {
        $addedIdx = $i
    }
}

if ($subSubIdx -eq -1) { throw "Could not find paragraph 'This is a sub-sub test'" }
if ($addedIdx -eq -1) { throw "Could not find paragraph 'Added this new change'" }
if ($addedIdx -le $subSubIdx) { throw "Unexpected paragraph order" }

# --- Step 1: rewrite "This is a sub-sub test" as four runs + the bookmark -
# Only the inner (text) range is touched, so the <w:p> itself (and its
# paraId/rsid bookkeeping) is left completely alone.
$pSubSub = $d.Paragraphs($subSubIdx)
$rSubSub = $d.Range($pSubSub.Range.Start, $pSubSub.Range.End - 1)

$subSubBody = "<w:p>" + `
    "<w:r><w:t>This is a sub-</w:t></w:r>" + `
    "<w:r><w:t>sub</w:t></w:r>" + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    "<w:r><w:t>change</w:t></w:r>" + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
    "</w:p>"

$rSubSub.InsertXML((New-PackageXml $subSubBody))

# --- Step 2: drop the old bookmark from "Added this new change" ----------
# A single-paragraph InsertXML that doesn't re-mention a bookmark anchored
# inside it gets the bookmark preserved (re-wrapped around the new
# content) by this engine instead of removed. Spanning the replacement
# across a paragraph-mark boundary (i.e. rewriting it together with its
# immediately preceding paragraph) avoids that and actually drops it, so
# we always pull in one extra paragraph of left-hand context here.
$pAdded = $d.Paragraphs($addedIdx)
$pPrevIdx = $addedIdx - 1
$pPrev = $d.Paragraphs($pPrevIdx)

$prevPPr = Get-ParagraphPPr $pPrev
$prevText = Get-ParagraphText $pPrev
$prevBody = "<w:p>$prevPPr" + (Get-XmlEscapedRun $prevText) + "</w:p>"

$addedPPr = Get-ParagraphPPr $pAdded
$addedBody = "<w:p>$addedPPr" + `
    '<w:r><w:t xml:space="preserve">Added this new </w:t></w:r>' + `
    "<w:r><w:t>change</w:t></w:r>" + `
    "</w:p>"

$rSpan = $d.Range($pPrev.Range.Start, $pAdded.Range.End - 1)
$rSpan.InsertXML((New-PackageXml ($prevBody + $addedBody)))

Write-Output "Edited paragraph text: '$(Get-ParagraphText $d.Paragraphs($subSubIdx))'"
Write-Output "Bookmark-host paragraph text: '$(Get-ParagraphText $d.Paragraphs($addedIdx))'"
